$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "The Future of AI in Education"

$ws.Range("B2").Value = "• Artificial Intelligence in Education (AIEd) - The use of AI technologies to enhance learning experiences`n• Personalized Learning - Tailoring educational content to individual student needs`n• Learning Analytics - The measurement and analysis of learning data to improve outcomes"

$ws.Range("C2").Value = "• Directly relevant to Livia's interests in leveraging AI for educational equity`n• Connects to her work with marginalized communities and learning design`n• Provides insights into career readiness and K-12 education applications"
